$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.046.28"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.208.73"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.44"
$ws.Range("E5").Value = "  +4.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.08"
$ws.Range("E6").Value = "  +2.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.208.34"
$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.68"
$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.736.44"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.218.82"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("E17").Value = "  +3.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.213.58"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.61"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("E21").Value = "  +5.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.02"
$ws.Range("E23").Value = "  +2.77%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.21"
$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.17"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.02"
$ws.Range("E27").Value = "  +3.17%  "

$ws.Range("E28").Value = "  +2.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +2.64%  "

$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.88"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("E31").Value = "  +9.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.15"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.35"
$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "484.62"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -4.62%  "

$ws.Range("E41").Value = "  +2.25%  "

$ws.Range("E42").Value = "  +3.89%  "

$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("E44").Value = "  +5.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.960.48"
$ws.Range("E45").Value = "  -3.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0644"
$ws.Range("E46").Value = "  +5.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.85"
$ws.Range("E47").Value = "  -0.82%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.24"
$ws.Range("E51").Value = "  +5.89%  "
